# Update the Metadata sheet: Version, Date, Contact
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item(1)
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# Add the new worksheet "Include from FSIII 2" as a copy-like sibling of
# "Include from FSIII", placed after it (at the end of the workbook).
$src = $wb.Worksheets.Item("Include from FSIII")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $last)
$new.Name = "Include from FSIII 2"

$new.Range("A1").Value = "Property"
$new.Range("B1").Value = "Operation"
$new.Range("C1").Value = "Value"
$new.Range("A2").Value = "concept"
$new.Range("B2").Value = "descendent-of"
$new.Range("C2").Value = "E"
$new.Range("A3").Value = "'"
$new.Range("B3").Value = "'"
$new.Range("A4").Value = "System URI"
$new.Range("B4").Value = "urn:oid:1.2.208.176.2.21"

# Match the formatting (borders/fills/fonts) of the source sheet.
$src.Range("A1:C1").Copy()
$new.Range("A1:C1").PasteSpecial(-4122)
$src.Range("A2:C2").Copy()
$new.Range("A2:C2").PasteSpecial(-4122)
$src.Range("A3:B3").Copy()
$new.Range("A3:B3").PasteSpecial(-4122)
$src.Range("A4:B4").Copy()
$new.Range("A4:B4").PasteSpecial(-4122)

# Match column widths.
$new.Columns.Item(1).ColumnWidth = 29.76
$new.Columns.Item(2).ColumnWidth = 49.76
